$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 29 (pushing existing rows 29-37 down to 31-39)
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(29).Insert()

# Fill in new row 29 (weekly record)
$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(29, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(29, 4).Value = 44855
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = 100112028
$ws.Cells.Item(29, 7).Value = "Sandia"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 300
$ws.Cells.Item(29, 11).Value = 750
$ws.Cells.Item(29, 12).Value = 800
$ws.Cells.Item(29, 13).Value = 775
$ws.Cells.Item(29, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(29, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(29, 16).Value = 775
$ws.Cells.Item(29, 17).Value = 1
$ws.Cells.Item(29, 18).Value = "Hortaliza"

# Fill in new row 30 (weekly record)
$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(30, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(30, 4).Value = 44855
$ws.Cells.Item(30, 5).Value = 15
$ws.Cells.Item(30, 6).Value = 100112028
$ws.Cells.Item(30, 7).Value = "Sandia"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Tercera"
$ws.Cells.Item(30, 10).Value = 340
$ws.Cells.Item(30, 11).Value = 750
$ws.Cells.Item(30, 12).Value = 800
$ws.Cells.Item(30, 13).Value = 775
$ws.Cells.Item(30, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(30, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(30, 16).Value = 775
$ws.Cells.Item(30, 17).Value = 1
$ws.Cells.Item(30, 18).Value = "Hortaliza"
